$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new FX observation row (row 5), mirroring the formatting of the
# existing USD/HKD rows above (e.g. row 2) so the date keeps its date format.
$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A5").Value = 44012        # 2020-06-30

$ws.Range("B5").Value = "USD"
$ws.Range("C5").Value = "HKD"
$ws.Range("D5").Value = 7.7504999999999997

# Leave the selection where the user finished entering data.
$ws.Range("D5").Select()
